$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.4738
$ws.Range("A3").Value = -21.9866
$ws.Range("A14").Value = -21.6816
$ws.Range("A16").Value = -21.61049999999999
$ws.Range("E18").Value = 18.24540000000002
$ws.Range("A21").Value = -19.90399999999998
$ws.Range("A23").Value = -20.28849999999998
$ws.Range("E24").Value = 16.4982
$ws.Range("A25").Value = -21.68889999999999
$ws.Range("E25").Value = 17.01649999999999
$ws.Range("A26").Value = -21.03239999999996
$ws.Range("E27").Value = 16.86799999999998
$ws.Range("A29").Value = -20.67699999999997
$ws.Range("E30").Value = 15.5193
$ws.Range("E31").Value = 16.07799999999999
$ws.Range("E39").Value = 16.0467
$ws.Range("A40").Value = -20.0252
$ws.Range("E42").Value = 16.52249999999999
$ws.Range("E48").Value = 17.39439999999999
$ws.Range("E51").Value = 17.5078
$ws.Range("E52").Value = 17.014
$ws.Range("A53").Value = -21.8976
$ws.Range("E55").Value = 16.6081
$ws.Range("E56").Value = 16.72079999999999
$ws.Range("A57").Value = -22.69180000000003
$ws.Range("E57").Value = 16.72500000000001
$ws.Range("A59").Value = -22.4299
$ws.Range("E60").Value = 15.6659
$ws.Range("A65").Value = -21.79319999999998
$ws.Range("A69").Value = -21.64139999999999
$ws.Range("E73").Value = 17.46880000000001
$ws.Range("E74").Value = 17.05149999999998
$ws.Range("A79").Value = -20.35450000000001
$ws.Range("A83").Value = -21.76619999999999
$ws.Range("E89").Value = 17.29550000000001
$ws.Range("E90").Value = 16.58499999999999
$ws.Range("A91").Value = -21.51220000000002
$ws.Range("E92").Value = 18.65690000000003
$ws.Range("A93").Value = -20.77889999999999
$ws.Range("A100").Value = -21.76859999999999
